# "before demo on create entity&sendnotify"
#
# Refresh the sample phone-number values on row 2 of every sheet so the
# demo doesn't show stale numbers from the last run. The four phone
# columns are the same across all four sheets:
#   G2  -> MobileNumber
#   AF2 -> Enquiry_PhoneNumber
#   AV2 -> Lead_PN
#   AZ2 -> Sales_PN

$wb = $excel.ActiveWorkbook

$newValues = @{
    "G2"  = "9840085196"
    "AF2" = "9840049097"
    "AV2" = "9840084872"
    "AZ2" = "9840010201"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in $newValues.Keys) {
        $cell = $ws.Range($addr)

        # A plain Value/Value2 assignment of a digit-only string gets
        # auto-converted to a number by Excel's "smart" cell-entry logic
        # (and an apostrophe-prefix / NumberFormat="@" workaround mutates
        # the cell style). Entering it as a text-producing formula and
        # then collapsing the formula to its literal value with a
        # values-only paste keeps both the text type AND the original
        # cell style/format untouched - exactly like the source data.
        $cell.Formula = '="' + $newValues[$addr] + '"'
        $cell.Copy() | Out-Null
        $cell.PasteSpecial(-4163) | Out-Null
    }
}

$excel.CutCopyMode = 0
